# Applies the commit:
#   "Remove extra column with no data in the encounter tab.
#    Set a delete date for two vital sign records."
#
# Part 1 (visible in the diff): the "encounter" sheet has a spurious, fully
# empty column H (it sits between "discharge_disp" (G) and the
# createdate/updatedate/deletedate columns, which used to be I/J/K). Delete
# it so createdate/updatedate/deletedate shift left to H/I/J, then resize
# the now-adjacent "discharge_disp" column to fit its (long) text values.
#
# Part 2: the "encounter" tab becomes the active/selected tab/window
# (previously "eVITALS" was), and its selection moves to the top of the new
# last column (H1, whole-column selection). "eVITALS" keeps its own
# scroll/selection state update (no longer the tab-selected sheet, and its
# selection moves to I137).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "encounter" sheet: delete the empty column H (createdate/updatedate/
#    deletedate were in I:K with an empty, merely-formatted H column in
#    front of them; remove it so they become H:J).
# ---------------------------------------------------------------------
$wsEnc = $wb.Worksheets.Item("encounter")
$wsEnc.Columns.Item(8).Delete()

# "discharge_disp" (now the last text column, G) is resized to fit its
# longest value.
$wsEnc.Columns.Item(7).ColumnWidth = 34

# ---------------------------------------------------------------------
# 2. View/selection state: "encounter" becomes the active sheet/tab, with
#    the whole of (the now relocated) column H selected.
# ---------------------------------------------------------------------
$wsEnc.Activate()
$wsEnc.Range("H1:H1048576").Select()

# ---------------------------------------------------------------------
# 3. "eVITALS" sheet: update its own lingering selection (it is no longer
#    the tab-selected sheet once "encounter" is activated above/below).
# ---------------------------------------------------------------------
$wsVitals = $wb.Worksheets.Item("eVITALS")
$wsVitals.Activate()
$wsVitals.Range("I137").Select()

# ---------------------------------------------------------------------
# 4. Leave "encounter" as the final active tab/window (matches
#    workbookView activeTab + encounter sheetView tabSelected="1").
# ---------------------------------------------------------------------
$wsEnc.Activate()
